$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-100 down to 35-101
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new data record
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 44662
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 100112030
$ws.Range("G34").Value = "Poroto granado"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 100
$ws.Range("K34").Value = 19000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = 19500
$ws.Range("N34").Value = "$/saco 25 kilos"
$ws.Range("O34").Value = "Provincia de Diguillín"
$ws.Range("P34").Value = 780
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
